$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New log entry (row 10): 2014-09-10, 20:00 -> 22:20, 0 min interruption, phase 5
$ws.Range("A10").Value = 41892
$ws.Range("B10").Value = 0.83333333333333337
$ws.Range("C10").Value = 0.93055555555555547
$ws.Range("D10").Value = 0
$ws.Range("E10").Formula = "=((HOUR(C10)-HOUR(B10))*60)+(MINUTE(C10)-MINUTE(B10))-D10"
$ws.Range("F10").Value = 5
$ws.Range("H10").Value = "Participé en el diagrama de caso de uso y escenarios de atributos de calidad"

# Wrapped comment text needs the taller row, same as rows 8/9
$ws.Rows.Item(10).RowHeight = 30

# Move/save the current selection like the author left it
$ws.Range("H16").Select()

Write-Output "done"
